# Weekly update: insert a new Jengibre price record at row 100, pushing the
# existing rows 100-199 down to 101-200 (dimension grows from A1:R199 to
# A1:R200). Row 99 and earlier stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 100; this shifts rows 100..199 -> 101..200
# and carries the existing row formatting (style s="2" on column D) down
# with them / into the new row.
$ws.Rows(100).Insert()

# Populate the newly inserted row 100 with this week's record. The
# non-varying columns (A,B,C,E,F,G,H,I,N,O,Q,R) match every other row in
# this Jengibre / Terminal La Palmera de La Serena subset.
$ws.Range("A100").Value = 8
$ws.Range("B100").Value = "Terminal La Palmera de La Serena"
$ws.Range("C100").Value = "Coquimbo"
$ws.Range("D100").Value = 45240
$ws.Range("E100").Value = 4
$ws.Range("F100").Value = 100114007
$ws.Range("G100").Value = "Jengibre"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 320
$ws.Range("K100").Value = 24000
$ws.Range("L100").Value = 25000
$ws.Range("M100").Value = 24500
$ws.Range("N100").Value = "`$/caja 13 kilos"
$ws.Range("O100").Value = "Perú"
$ws.Range("P100").Value = 1885
$ws.Range("Q100").Value = 13
$ws.Range("R100").Value = "Hortaliza"
